# term/ValueSet-KLInterventionsFSIII.xlsx : bump published version 1.0.0 -> 1.1.0
# and refresh the "Date" metadata row to match the new release.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 1).Text

    if ($label -eq "Version") {
        $ws.Cells.Item($r, 2).Value = "1.1.0"
    }
    elseif ($label -eq "Date") {
        $ws.Cells.Item($r, 2).Value = "2023-07-10T23:08:03+02:00"
    }
}
